$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text to uppercase
$ws.Range("A1").Value = "EMERGENCIAS PREVIAMENTE CARGADAS"

# Apply Excel's built-in "Bad" cell style (Incorrecto) and bold it
$ws.Range("A1").Style = "Bad"
$ws.Range("A1").Font.Bold = $true

# Widen column A to fit new text
$ws.Columns.Item(1).ColumnWidth = 38
